$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new data rows at 1109-1111 (new survey date 2022-05-25 / serial 44706).
# This shifts the existing rows 1109-1167 down to 1112-1170, which automatically
# reproduces the tail rows 1168-1170 (former 1165-1167) with no further edits needed.
$ws.Range("A1109:A1111").EntireRow.Insert()

# Row 1109
$ws.Cells.Item(1109, 1).Value = 9
$ws.Cells.Item(1109, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(1109, 3).Value = 'Metropolitana'
$ws.Cells.Item(1109, 4).Value = 44706
$ws.Cells.Item(1109, 5).Value = 13
$ws.Cells.Item(1109, 6).Value = 100114001
$ws.Cells.Item(1109, 7).Value = 'Papa'
$ws.Cells.Item(1109, 8).Value = 'Asterix'
$ws.Cells.Item(1109, 9).Value = '1a nueva(o)'
$ws.Cells.Item(1109, 10).Value = 170
$ws.Cells.Item(1109, 11).Value = 9000
$ws.Cells.Item(1109, 12).Value = 10000
$ws.Cells.Item(1109, 13).Value = 9412
$ws.Cells.Item(1109, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(1109, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(1109, 16).Value = 376
$ws.Cells.Item(1109, 17).Value = 25
$ws.Cells.Item(1109, 18).Value = 'Hortaliza'

# Row 1110
$ws.Cells.Item(1110, 1).Value = 9
$ws.Cells.Item(1110, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(1110, 3).Value = 'Metropolitana'
$ws.Cells.Item(1110, 4).Value = 44706
$ws.Cells.Item(1110, 5).Value = 13
$ws.Cells.Item(1110, 6).Value = 100114001
$ws.Cells.Item(1110, 7).Value = 'Papa'
$ws.Cells.Item(1110, 8).Value = 'Rodeo'
$ws.Cells.Item(1110, 9).Value = '1a (cosecha lavada)'
$ws.Cells.Item(1110, 10).Value = 240
$ws.Cells.Item(1110, 11).Value = 9000
$ws.Cells.Item(1110, 12).Value = 10000
$ws.Cells.Item(1110, 13).Value = 9458
$ws.Cells.Item(1110, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(1110, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(1110, 16).Value = 378
$ws.Cells.Item(1110, 17).Value = 25
$ws.Cells.Item(1110, 18).Value = 'Hortaliza'

# Row 1111
$ws.Cells.Item(1111, 1).Value = 9
$ws.Cells.Item(1111, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(1111, 3).Value = 'Metropolitana'
$ws.Cells.Item(1111, 4).Value = 44706
$ws.Cells.Item(1111, 5).Value = 13
$ws.Cells.Item(1111, 6).Value = 100114001
$ws.Cells.Item(1111, 7).Value = 'Papa'
$ws.Cells.Item(1111, 8).Value = 'Rodeo'
$ws.Cells.Item(1111, 9).Value = '1a (cosecha)'
$ws.Cells.Item(1111, 10).Value = 190
$ws.Cells.Item(1111, 11).Value = 8000
$ws.Cells.Item(1111, 12).Value = 8000
$ws.Cells.Item(1111, 13).Value = 8000
$ws.Cells.Item(1111, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(1111, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(1111, 16).Value = 320
$ws.Cells.Item(1111, 17).Value = 25
$ws.Cells.Item(1111, 18).Value = 'Hortaliza'
